$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for "Perejil" (Vega Central Mapocho de
# Santiago). Insert a new row above the current row 418, shifting the
# existing rows 418-443 down to 419-444, then fill in the new record's data.
$ws.Rows(418).Insert()

$ws.Range("A418").Value = 9
$ws.Range("B418").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C418").Value = "Metropolitana"
$ws.Range("D418").Value = 44826
$ws.Range("E418").Value = 13
$ws.Range("F418").Value = 100112044
$ws.Range("G418").Value = "Perejil"
$ws.Range("H418").Value = "Sin especificar"
$ws.Range("I418").Value = "Primera"
$ws.Range("J418").Value = 70
$ws.Range("K418").Value = 10000
$ws.Range("L418").Value = 10000
$ws.Range("M418").Value = 10000
$ws.Range("N418").Value = "$/docena de atados"
$ws.Range("O418").Value = "Región Metropolitana"
$ws.Range("P418").Value = 3333
$ws.Range("Q418").Value = 3
$ws.Range("R418").Value = "Hortaliza"
